$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) for columns B:E -> 15, 16, 15, 16
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: B2 and D2 get new values; C2 and E2 are cleared (no longer populated)
$ws.Range("B2").Value = 2.3772960904422913
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 3.479341849987208
$ws.Range("E2").ClearContents()

# Row 3: B3:E3 all get new values
$ws.Range("B3").Value = 2.2032100745536449
$ws.Range("C3").Value = -0.75226107008933984
$ws.Range("D3").Value = 3.6460566198073323
$ws.Range("E3").Value = -0.097459496988475572

# Update the sheet selection to match the narrower post-edit range
$ws.Range("B1:E3").Select()
